$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be stored as text even if the value looks numeric,
    # then restore the default "Normal" style so no stray number format sticks.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 14 / Row 15: Polygon and Polkadot swap ranking position, with updated data
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D14") "5.106"
$ws.Range("E14").Value = "  +4.43%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D15") "0.6826"
$ws.Range("E15").Value = "  +2.18%  "

# Remaining rows: refreshed price (D) and 1h volume change (E) values
$ws.Range("D2").Value = "30.738.87"
$ws.Range("E2").Value = "  +2.60%  "
$ws.Range("D3").Value = "1.894.23"
$ws.Range("E3").Value = "  +0.79%  "
Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  +0.16%  "
Set-TextValue $ws.Range("D5") "246.91"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("E6").Value = "  +0.19%  "
Set-TextValue $ws.Range("D7") "0.4933"
$ws.Range("E7").Value = "  -0.83%  "
Set-TextValue $ws.Range("D8") "0.2954"
$ws.Range("E8").Value = "  +1.01%  "
Set-TextValue $ws.Range("D9") "0.06816"
$ws.Range("E9").Value = "  +2.77%  "
$ws.Range("D10").Value = "1.894.21"
$ws.Range("E10").Value = "  +0.80%  "
Set-TextValue $ws.Range("D11") "17.31"
$ws.Range("E11").Value = "  +3.39%  "
Set-TextValue $ws.Range("D12") "92.32"
$ws.Range("E12").Value = "  +6.91%  "
Set-TextValue $ws.Range("D13") "0.07251"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D16").Value = "30.706.36"
$ws.Range("E16").Value = "  +2.53%  "
Set-TextValue $ws.Range("D17") "0.000008009"
$ws.Range("E17").Value = "  +1.35%  "
Set-TextValue $ws.Range("D18") "13.33"
$ws.Range("E18").Value = "  +4.32%  "
Set-TextValue $ws.Range("D19") "1.001"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").Value = "2.138.47"
$ws.Range("E20").Value = "  +0.88%  "
Set-TextValue $ws.Range("D21") "1.002"
$ws.Range("E21").Value = "  +0.40%  "
Set-TextValue $ws.Range("D22") "4.855"
$ws.Range("E22").Value = "  +1.73%  "
Set-TextValue $ws.Range("D23") "193.13"
$ws.Range("E23").Value = "  +36.03%  "
Set-TextValue $ws.Range("D24") "6.060"
$ws.Range("E24").Value = "  +7.15%  "
Set-TextValue $ws.Range("D25") "9.417"
$ws.Range("E25").Value = "  +3.92%  "
Set-TextValue $ws.Range("D26") "155.70"
$ws.Range("E26").Value = "  +4.27%  "
Set-TextValue $ws.Range("D27") "19.23"
$ws.Range("E27").Value = "  +12.05%  "
$ws.Range("E28").Value = "  +0.81%  "
Set-TextValue $ws.Range("D29") "1.401"
$ws.Range("E29").Value = "  +0.80%  "
Set-TextValue $ws.Range("D30") "4.370"
$ws.Range("E30").Value = "  +4.55%  "
$ws.Range("E31").Value = "  +2.67%  "
$ws.Range("E32").Value = "  +2.59%  "
Set-TextValue $ws.Range("D33") "0.05212"
$ws.Range("E33").Value = "  +2.87%  "
Set-TextValue $ws.Range("D34") "0.7464"
$ws.Range("E34").Value = "  +5.26%  "
Set-TextValue $ws.Range("D35") "1.127"
$ws.Range("E35").Value = "  +1.88%  "
Set-TextValue $ws.Range("D36") "2.722"
$ws.Range("E36").Value = "  +2.14%  "
Set-TextValue $ws.Range("D37") "0.01864"
$ws.Range("E37").Value = "  +5.00%  "
$ws.Range("E38").Value = "  -0.35%  "
Set-TextValue $ws.Range("D39") "2.167"
$ws.Range("E39").Value = "  -0.46%  "
Set-TextValue $ws.Range("D40") "0.9377"
$ws.Range("E40").Value = "  +0.80%  "
Set-TextValue $ws.Range("D41") "0.4448"
$ws.Range("E41").Value = "  +4.60%  "
Set-TextValue $ws.Range("D42") "106.64"
$ws.Range("E42").Value = "  +4.21%  "
Set-TextValue $ws.Range("D43") "5.787"
$ws.Range("E43").Value = "  -0.23%  "
Set-TextValue $ws.Range("D45") "7.705"
$ws.Range("E45").Value = "  +2.97%  "
Set-TextValue $ws.Range("D46") "0.1346"
$ws.Range("E46").Value = "  +6.86%  "
Set-TextValue $ws.Range("D47") "0.05860"
$ws.Range("E47").Value = "  +3.69%  "
Set-TextValue $ws.Range("D48") "8.778"
$ws.Range("E48").Value = "  +6.42%  "
Set-TextValue $ws.Range("D49") "1.443"
$ws.Range("E49").Value = "  +8.35%  "
Set-TextValue $ws.Range("D50") "0.3955"
$ws.Range("E50").Value = "  +5.03%  "
Set-TextValue $ws.Range("D51") "33.68"
$ws.Range("E51").Value = "  +3.88%  "
